$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2580.8157
$ws.Range("I135").Value = 2740.862
$ws.Range("J135").Value = 2065.111
$ws.Range("K135").Value = 24667.758
$ws.Range("L135").Value = 18585.999
$ws.Range("M135").Value = -22132.758
$ws.Range("N135").Value = -23655.999

$ws.Range("H137").Value = 1843.8948
$ws.Range("J137").Value = 2536.6316
$ws.Range("L137").Value = 7609.8948
$ws.Range("N137").Value = -12709.8948

$ws.Range("H138").Value = 2707.295
$ws.Range("I138").Value = 1147.2
$ws.Range("J138").Value = 4349.5
$ws.Range("K138").Value = 3441.6
$ws.Range("L138").Value = 13048.5
$ws.Range("M138").Value = 1698.4
$ws.Range("N138").Value = -23328.5

$ws.Range("H141").Value = 1598.8667
$ws.Range("I141").Value = 1598.8667
$ws.Range("K141").Value = 4796.6001
$ws.Range("M141").Value = 383.3999000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3979.763
$ws.Range("I32").Value = 3636.3918
$ws.Range("J32").Value = 5084.522
$ws.Range("K32").Value = 3636.3918
$ws.Range("L32").Value = 5084.522
$ws.Range("M32").Value = -3349.3918
$ws.Range("N32").Value = -5658.522

$ws.Range("H45").Value = 7985.9546
$ws.Range("I45").Value = 9775.706
$ws.Range("K45").Value = 9775.706
$ws.Range("M45").Value = -9398.706

$ws.Range("H61").Value = 3336.83
$ws.Range("I61").Value = 4126.8203
$ws.Range("J61").Value = 1136.1428
$ws.Range("K61").Value = 4126.8203
$ws.Range("L61").Value = 1136.1428
$ws.Range("M61").Value = -3914.8203
$ws.Range("N61").Value = -1560.1428

$ws.Range("H132").Value = 4653.273
$ws.Range("I132").Value = 3750
$ws.Range("J132").Value = 5169.4287
$ws.Range("K132").Value = 11250
$ws.Range("L132").Value = 15508.2861
$ws.Range("M132").Value = -8720
$ws.Range("N132").Value = -20568.2861

$ws.Range("H136").Value = 3336.83
$ws.Range("I136").Value = 4126.8203
$ws.Range("J136").Value = 1136.1428
$ws.Range("K136").Value = 12380.4609
$ws.Range("L136").Value = 3408.4284
$ws.Range("M136").Value = -9830.460900000002
$ws.Range("N136").Value = -8508.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 8759.799999999999
$ws.Range("I54").Value = 6449.75
$ws.Range("J54").Value = 18000
$ws.Range("K54").Value = 6449.75
$ws.Range("L54").Value = 18000
$ws.Range("M54").Value = -5965.75
$ws.Range("N54").Value = -18968

$ws.Range("H63").Value = 29000
$ws.Range("J63").Value = 29000
$ws.Range("L63").Value = 29000
$ws.Range("N63").Value = -30372

$ws.Range("H66").Value = 29000
$ws.Range("J66").Value = 29000
$ws.Range("L66").Value = 87000
$ws.Range("N66").Value = -93864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3605.7073
$ws.Range("I31").Value = 1763.9565
$ws.Range("J31").Value = 4323.6777
$ws.Range("K31").Value = 1763.9565
$ws.Range("L31").Value = 4323.6777
$ws.Range("M31").Value = -1468.9565
$ws.Range("N31").Value = -4913.6777

$ws.Range("H34").Value = 3605.7073
$ws.Range("I34").Value = 1763.9565
$ws.Range("J34").Value = 4323.6777
$ws.Range("K34").Value = 1763.9565
$ws.Range("L34").Value = 4323.6777
$ws.Range("M34").Value = -1561.9565
$ws.Range("N34").Value = -4727.6777

$ws.Range("H122").Value = 908.76
$ws.Range("I122").Value = 912.9375
$ws.Range("J122").Value = 901.3333
$ws.Range("K122").Value = 2738.8125
$ws.Range("L122").Value = 2703.9999
$ws.Range("M122").Value = -288.8125
$ws.Range("N122").Value = -7603.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 794.4286
$ws.Range("I45").Value = 428
$ws.Range("J45").Value = 855.5
$ws.Range("K45").Value = 1284
$ws.Range("L45").Value = 2566.5
$ws.Range("M45").Value = -752
$ws.Range("N45").Value = -3630.5

$ws.Range("H68").Value = 2680.4
$ws.Range("I68").Value = 4073.5
$ws.Range("J68").Value = 1751.6666
$ws.Range("K68").Value = 12220.5
$ws.Range("L68").Value = 5254.9998
$ws.Range("M68").Value = -11409.5
$ws.Range("N68").Value = -6876.9998

$ws.Range("H71").Value = 2680.4
$ws.Range("I71").Value = 4073.5
$ws.Range("J71").Value = 1751.6666
$ws.Range("K71").Value = 36661.5
$ws.Range("L71").Value = 15764.9994
$ws.Range("M71").Value = -32605.5
$ws.Range("N71").Value = -23876.9994

$ws.Range("H96").Value = 6345.2
$ws.Range("J96").Value = 7425
$ws.Range("L96").Value = 22275
$ws.Range("N96").Value = -26393

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 751.7059
$ws.Range("I97").Value = 790
$ws.Range("J97").Value = 697
$ws.Range("K97").Value = 790
$ws.Range("L97").Value = 697
$ws.Range("M97").Value = -294
$ws.Range("N97").Value = -1689

$ws.Range("H132").Value = 4612.8076
$ws.Range("I132").Value = 9800
$ws.Range("J132").Value = 3669.682
$ws.Range("K132").Value = 29400
$ws.Range("L132").Value = 11009.046
$ws.Range("M132").Value = -26870
$ws.Range("N132").Value = -16069.046

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5823198
$ws.Range("I122").Value = 7147776.5
$ws.Range("K122").Value = 21443329.5
$ws.Range("M122").Value = -21440879.5

$ws.Range("H132").Value = 22816118
$ws.Range("I132").Value = 28899084
$ws.Range("J132").Value = 4998.25
$ws.Range("K132").Value = 86697252
$ws.Range("L132").Value = 14994.75
$ws.Range("M132").Value = -86694722
$ws.Range("N132").Value = -20054.75

$ws.Range("H136").Value = 5382.375
$ws.Range("I136").Value = 4441.6
$ws.Range("J136").Value = 7734.3125
$ws.Range("K136").Value = 13324.8
$ws.Range("L136").Value = 23202.9375
$ws.Range("M136").Value = -10774.8
$ws.Range("N136").Value = -28302.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 24125
$ws.Range("J86").Value = 24125
$ws.Range("L86").Value = 24125
$ws.Range("N86").Value = -26371

$ws.Range("H89").Value = 24125
$ws.Range("J89").Value = 24125
$ws.Range("L89").Value = 120625
$ws.Range("N89").Value = -131857

$ws.Range("H132").Value = 1616.1628
$ws.Range("I132").Value = 1467.6333
$ws.Range("J132").Value = 1958.9231
$ws.Range("K132").Value = 4402.8999
$ws.Range("L132").Value = 5876.7693
$ws.Range("M132").Value = -1872.8999
$ws.Range("N132").Value = -10936.7693

$ws.Range("H136").Value = 3159.0322
$ws.Range("I136").Value = 3692.4443
$ws.Range("J136").Value = 2420.4614
$ws.Range("K136").Value = 11077.3329
$ws.Range("L136").Value = 7261.3842
$ws.Range("M136").Value = -8527.332900000001
$ws.Range("N136").Value = -12361.3842
